$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (2023-09-05 -> 2023-09-06) for rows 2 through 9.
$ws.Range("C2:C9").Value = 45175
